$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Merge the "I, {{RANK}} {{NAME}} {{BADGE}} ... just, probable and
#    reasonable cause ..." sentence back into a single run (removing
#    the bold "probable" run + its surrounding proofErr markers).
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "there is just, probable and reasonable cause",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "there is just, probable and reasonable cause", 2) | Out-Null

# -----------------------------------------------------------------
# 2) Split every "Judge, {{COURT}}" run into three runs:
#      "Judge"  +  " {{JUDGE}}"  +  ", {{COURT}}"
#    Each occurrence keeps identical run formatting (sz/szCs 24), so
#    a plain InsertAfter would normally get re-coalesced back into a
#    single run. Briefly nudging the new run's font size to a
#    different value keeps it a distinct run; resetting that size to
#    match (12pt == sz 24) in a later, independent Find pass leaves
#    the run boundaries alone while restoring identical formatting.
# -----------------------------------------------------------------
$searchFrom = 0
while ($true) {
    $rng = $d.Range($searchFrom, $d.Content.End)
    $found = $rng.Find.Execute(
        "Judge, {{COURT}}",
        $false, $false, $false, $false, $false, $true, 1, $false,
        "", 0)
    if (-not $found) { break }

    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # Collapse to the point right after "Judge" and insert the new
    # merge-field text there.
    $insPoint = $d.Range($matchStart + 5, $matchStart + 5)
    $insPoint.InsertAfter(" {{JUDGE}}")

    # Re-find the freshly inserted text so we can nudge its size.
    $judgeRng = $d.Range($matchStart, $matchEnd + 10)
    $judgeRng.Find.Execute(
        " {{JUDGE}}", $false, $false, $false, $false, $false, $true,
        1, $false, "", 0) | Out-Null
    $judgeRng.Font.Size = 11

    $searchFrom = $matchEnd + 10
}

# Now restore the nudged runs' font size to match their neighbours
# (12pt -> sz 24) without re-merging the run boundaries.
$searchFrom = 0
while ($true) {
    $rng = $d.Range($searchFrom, $d.Content.End)
    $found = $rng.Find.Execute(
        " {{JUDGE}}", $false, $false, $false, $false, $false, $true,
        1, $false, "", 0)
    if (-not $found) { break }
    $rng.Font.Size = 12
    $searchFrom = $rng.End
}
